# Update cryptos list with latest scraped price/volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    # Force text storage so numeric-looking strings (e.g. "247.01")
    # are not reinterpreted as numbers, then drop back to the default
    # (unstyled) cell style so no stray formatting is introduced.
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '35.515.81'
Set-TextCell 'E2' '  +1.30%  '
Set-TextCell 'D3' '1.909.38'
Set-TextCell 'E3' '  +2.90%  '
Set-TextCell 'E4' '  +0.23%  '
Set-TextCell 'D5' '247.01'
Set-TextCell 'E5' '  +4.16%  '
Set-TextCell 'E6' '  +5.49%  '
Set-TextCell 'E7' '  +0.26%  '
Set-TextCell 'D8' '42.15'
Set-TextCell 'E8' '  +0.15%  '
Set-TextCell 'E9' '  +5.48%  '
Set-TextCell 'D10' '49.02'
Set-TextCell 'E10' '  +5.09%  '
Set-TextCell 'D11' '0.0717'
Set-TextCell 'E11' '  +3.22%  '
Set-TextCell 'E12' '  +1.01%  '
Set-TextCell 'D13' '2.189.09'
Set-TextCell 'E13' '  +3.02%  '
Set-TextCell 'D14' '12.32'
Set-TextCell 'E14' '  +8.04%  '
Set-TextCell 'D15' '0.700'
Set-TextCell 'E15' '  +3.54%  '
Set-TextCell 'B16' 'WrappedEther'
Set-TextCell 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D16' '1.901.55'
Set-TextCell 'E16' '  +2.37%  '
Set-TextCell 'B17' 'Polkadot'
Set-TextCell 'C17' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D17' '4.86'
Set-TextCell 'E17' '  +3.66%  '
Set-TextCell 'D18' '35.536.57'
Set-TextCell 'E18' '  +1.44%  '
Set-TextCell 'D19' '72.38'
Set-TextCell 'E19' '  +2.95%  '
Set-TextCell 'D20' '0.0₃0828'
Set-TextCell 'E20' '  +4.35%  '
Set-TextCell 'D21' '244.62'
Set-TextCell 'E21' '  +1.76%  '
Set-TextCell 'D22' '12.73'
Set-TextCell 'E22' '  +5.09%  '
Set-TextCell 'D23' '4.84'
Set-TextCell 'E23' '  +2.08%  '
Set-TextCell 'E25' '  +1.47%  '
Set-TextCell 'B26' 'Monero'
Set-TextCell 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D26' '171.69'
Set-TextCell 'E26' '  +0.23%  '
Set-TextCell 'B27' 'PancakeSwap'
Set-TextCell 'C27' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 'D27' '2.20'
Set-TextCell 'E27' '  +16.78%  '
Set-TextCell 'D28' '8.52'
Set-TextCell 'E28' '  +7.26%  '
Set-TextCell 'E29' '  +3.96%  '
Set-TextCell 'E30' '  +3.92%  '
Set-TextCell 'D31' '0.972'
Set-TextCell 'E31' '  +24.05%  '
Set-TextCell 'D32' '4.16'
Set-TextCell 'E32' '  +4.18%  '
Set-TextCell 'D33' '0.0572'
Set-TextCell 'E33' '  +2.69%  '
Set-TextCell 'E34' '  +5.19%  '
Set-TextCell 'E35' '  +0.18%  '
Set-TextCell 'D36' '1.74'
Set-TextCell 'E36' '  +6.33%  '
Set-TextCell 'E37' '  +1.17%  '
Set-TextCell 'E38' '  +2.90%  '
Set-TextCell 'D39' '1.11'
Set-TextCell 'E39' '  +3.24%  '
Set-TextCell 'D40' '0.0209'
Set-TextCell 'E40' '  +2.58%  '
Set-TextCell 'D41' '92.66'
Set-TextCell 'E41' '  +1.17%  '
Set-TextCell 'D42' '0.0628'
Set-TextCell 'E42' '  +13.58%  '
Set-TextCell 'E43' '  +5.54%  '
Set-TextCell 'D44' '1.351.39'
Set-TextCell 'E44' '  -0.03%  '
Set-TextCell 'D45' '2.40'
Set-TextCell 'E45' '  +2.39%  '
Set-TextCell 'B46' 'Gas'
Set-TextCell 'C46' 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
Set-TextCell 'D46' '13.10'
Set-TextCell 'E46' '  +3.92%  '
Set-TextCell 'B47' 'MultiversX'
Set-TextCell 'C47' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell 'D47' '47.24'
Set-TextCell 'E47' '  +38.16%  '
Set-TextCell 'D48' '2.42'
Set-TextCell 'E48' '  +0.50%  '
Set-TextCell 'D49' '2.78'
Set-TextCell 'E49' '  +1.87%  '
Set-TextCell 'D50' '6.60'
Set-TextCell 'E50' '  +2.18%  '
Set-TextCell 'D51' '2.099.27'
Set-TextCell 'E51' '  +3.15%  '
